# Refresh the cryptos price/volume snapshot (rows 2-51, columns D & E).
# Values in column D that look like plain numbers (e.g. "1.001") are written
# with a leading apostrophe so Excel keeps them as text (matching the
# workbook's existing inline-string convention) instead of coercing them to
# numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.430.83'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = '1.827.75'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''314.70'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.5155'
$ws.Range("E7").Value = '  -2.49%  '
$ws.Range("D8").Value = '''0.3896'
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").Value = '''0.07635'
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("D11").Value = '''1.108'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").Value = '''21.06'
$ws.Range("E12").Value = '  +3.86%  '
$ws.Range("D13").Value = '''6.275'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '''7.536'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '1.824.55'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '''93.38'
$ws.Range("E17").Value = '  +5.58%  '
$ws.Range("D18").Value = '''0.00001083'
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").Value = '''0.06686'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").Value = '''17.65'
$ws.Range("E20").Value = '  +3.01%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '''6.177'
$ws.Range("E22").Value = '  +3.49%  '
$ws.Range("D23").Value = '28.456.03'
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").Value = '''11.15'
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("D25").Value = '''2.251'
$ws.Range("E25").Value = '  +7.59%  '
$ws.Range("D26").Value = '''156.92'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = '''20.61'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("D28").Value = '2.034.78'
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("D29").Value = '''2.396'
$ws.Range("E29").Value = '  +3.70%  '
$ws.Range("D30").Value = '''124.71'
$ws.Range("E30").Value = '  +2.19%  '
$ws.Range("D31").Value = '''1.117'
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("D32").Value = '''0.1086'
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").Value = '''5.662'
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = '''3.672'
$ws.Range("D35").Value = '''0.07037'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '''0.2221'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").Value = '''8.911'
$ws.Range("E37").Value = '  +6.37%  '
$ws.Range("D38").Value = '''0.02324'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").Value = '''5.121'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").Value = '''0.6299'
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("D41").Value = '''11.22'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''1.184'
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '''1.392'
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").Value = '''13.39'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '''0.5901'
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("D47").Value = '''3.711'
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").Value = '''124.58'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").Value = '''1.974'
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("D50").Value = '''1.201'
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("D51").Value = '''0.06925'
$ws.Range("E51").Value = '  +2.00%  '
